$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.131.63"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "3.123.01"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.21"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.30"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.120.76"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "3.659.32"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.65"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "58.173.84"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "3.128.98"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.34"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.66"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.39"
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.06"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.42"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.36"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E39").Value = "  -3.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.66"
$ws.Range("E40").Value = "  +14.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0667"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.693"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "3.162.76"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.61"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0262"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "2.283.41"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +4.06%  "
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.64"
$ws.Range("E51").Value = "  -0.17%  "
